$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix "1)Запустить тест " -> "1) Запустить тест " (missing space after ")")
#    in the two existing "Шаги" (Steps) cells F2 and F4.
# ---------------------------------------------------------------------------
$ws.Range("F2").Characters(1, 17).Text = "1) Запустить тест "
$ws.Range("F4").Characters(1, 17).Text = "1) Запустить тест "

# ---------------------------------------------------------------------------
# 2) Add two new test cases (rows 6-9 and 10-13), each spanning 4 rows with
#    merged cells for columns A, B, D, E, I (same layout as the existing two
#    test cases in rows 2-3 / 4-5).
# ---------------------------------------------------------------------------

# --- Test case 3 (rows 6-9) ------------------------------------------------
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Раздел News. Переход в раздел Control Panel. Проверка раздела Creating News"
$ws.Range("D6").Value = "Критическая"
$ws.Range("E6").Value = "ОС Windows 10. Android Studio Dolphin | 2021.3.1 Patch. AP I 29. Приложение на английском языке"
$ws.Range("F6").Value = "1) Перейти в раздел News."
$ws.Range("F7").Value = "2) Перейти в раздел Control Panel."
$ws.Range("F8").Value = "3) Перейти в раздел Creatin News"
$ws.Range("F9").Value = "4) Нажать на поле Category"
$ws.Range("G9").Value = "Открылось всплывающее меню с выбором категорий на английском языке."
$ws.Range("H9").Value = "Открылось всплывающее меню с выбором категорий на русском языке."
$ws.Hyperlinks.Add($ws.Range("I6"), "https://photos.app.goo.gl/QQFFuKMFYTnt4SsDA", "", "", "https://photos.app.goo.gl/QQFFuKMFYTnt4SsDA")

$ws.Range("A6:A9").Merge()
$ws.Range("B6:B9").Merge()
$ws.Range("D6:D9").Merge()
$ws.Range("E6:E9").Merge()
$ws.Range("I6:I9").Merge()

# --- Test case 4 (rows 10-13) ----------------------------------------------
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "Раздел News. Переход в раздел Control Panel. Проверка раздела Filter News."
$ws.Range("D10").Value = "Критическая"
$ws.Range("E10").Value = "ОС Windows 10. Android Studio Dolphin | 2021.3.1 Patch. AP I 29. Приложение на английском языке"
$ws.Range("F10").Value = "1) Перейти в раздел News."
$ws.Range("F11").Value = "2) Перейти в раздел Control Panel."
$ws.Range("F12").Value = "3) Перейти в раздел Filter News"
$ws.Range("F13").Value = "4) Нажать на поле Category"
$ws.Range("G13").Value = "Открылось всплывающее меню с выбором категорий на английском языке."
$ws.Range("H13").Value = "Открылось всплывающее меню с выбором категорий на русском языке."
$ws.Hyperlinks.Add($ws.Range("I10"), "https://photos.app.goo.gl/P1fFoRHGpX2a8v2LA", "", "", "https://photos.app.goo.gl/P1fFoRHGpX2a8v2LA")

$ws.Range("A10:A13").Merge()
$ws.Range("B10:B13").Merge()
$ws.Range("D10:D13").Merge()
$ws.Range("E10:E13").Merge()
$ws.Range("I10:I13").Merge()

# ---------------------------------------------------------------------------
# 3) Formatting: apply the same font/alignment used elsewhere in the sheet to
#    the newly added cells (XO Thames 12, centered + vertical centered,
#    wrapped for the "Шаги" column).
# ---------------------------------------------------------------------------
$newRange = $ws.Range("A6:I13")
$newRange.Font.Name = "XO Thames"
$newRange.Font.Size = 12
$newRange.HorizontalAlignment = -4108
$newRange.VerticalAlignment = -4108
$newRange.WrapText = $true
